# Feria Lagunitas de Puerto Montt - Poroto granado
# Weekly update: insert a new price-report row at row 16 (pushes the
# existing historical rows down by one, dimension grows from R27 to R28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 16, shifting rows 16-27 -> 17-28.
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with this week's market report.
$ws.Cells.Item(16, 1).Value  = 4
$ws.Cells.Item(16, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(16, 3).Value  = "Los Lagos"
$ws.Cells.Item(16, 4).Value  = 44589
$ws.Cells.Item(16, 5).Value  = 10
$ws.Cells.Item(16, 6).Value  = 100112030
$ws.Cells.Item(16, 7).Value  = "Poroto granado"
$ws.Cells.Item(16, 8).Value  = "Sin especificar"
$ws.Cells.Item(16, 9).Value  = "Primera"
$ws.Cells.Item(16, 10).Value = 90
$ws.Cells.Item(16, 11).Value = 31000
$ws.Cells.Item(16, 12).Value = 31000
$ws.Cells.Item(16, 13).Value = 31000
$ws.Cells.Item(16, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(16, 15).Value = "Región Metropolitana"
$ws.Cells.Item(16, 16).Value = 1240
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = "Hortaliza"
